$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list with latest price (column D) and 1h volume change (column E) values
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.127.06"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.68"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.56"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4574"
$ws.Range("E7").Value = "  +7.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3729"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07330"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8602"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.01"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.49"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.696"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.91"
$ws.Range("E14").Value = "  +5.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.344"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07076"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008836"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.126.51"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.188"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.006"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.60"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.227"
$ws.Range("E26").Value = "  +5.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.51"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.268"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.46"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08868"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7640"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.195"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.964"
$ws.Range("E33").Value = "  +5.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.469"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.105"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05289"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5356"
$ws.Range("E39").Value = "  +6.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.176"
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.891"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1715"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5228"
$ws.Range("E43").Value = "  +11.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.622"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.74"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.003"
$ws.Range("E46").Value = "  +11.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.07"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06477"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.0000"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9226"
$ws.Range("E51").Value = "  +1.26%  "
